# Add data for 2022-04-05: update sheet/labels for March through 03-28,
# and bump the March row + Total row figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab to reflect the new "through" date.
$ws.Name = "Through 2022-03-28"

# Update the label in A4 ("March (through 03-27)" -> "March (through 03-28)")
$ws.Range("A4").Value = "March (through 03-28)"

# Row 4 (March) updated counts
$ws.Range("B4").Value = 27
$ws.Range("D4").Value = 51
$ws.Range("E4").Value = 57
$ws.Range("F4").Value = 28
$ws.Range("G4").Value = 53
$ws.Range("H4").Value = 76
$ws.Range("I4").Value = 120

# Row 5 (Total) updated counts
$ws.Range("B5").Value = 64
$ws.Range("D5").Value = 182
$ws.Range("E5").Value = 194
$ws.Range("F5").Value = 107
$ws.Range("G5").Value = 194
$ws.Range("H5").Value = 418
$ws.Range("I5").Value = 420
